$d = $word.ActiveDocument

# 1. Split the paragraph about login into two paragraphs.
$d.Content.Find.Execute(
    "próprio site. O usuário pode traçar",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "próprio site.^pO usuário pode traçar",
    2)

# 2. Fix "investidos" -> "investido" + "r" (split into two runs effectively,
#    text result is "investidor").
$d.Content.Find.Execute(
    "Perfil de investidos, que",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Perfil de investidor, que",
    2)

# 3. Insert "5" before " minutos." in the Cotações paragraph (fix the double
#    space / missing number).
$d.Content.Find.Execute(
    "atualizada a cada  minutos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "atualizada a cada 5 minutos.",
    2)

# 4. Fix "postgresql" -> "PostgreSQL"
$d.Content.Find.Execute(
    "postgresql",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PostgreSQL",
    2)

# 5. Fix "messes" -> "meses"
$d.Content.Find.Execute(
    "10 messes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "10 meses.",
    2)
